# Update the "studyDesign" sheet with two new top rows capturing the
# study design name / description, and make that sheet the active tab
# (it was previously "studyDesignElements").

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("studyDesign")

# Insert two new blank rows at the very top of the sheet.
$ws.Rows.Item(1).Insert()
$ws.Rows.Item(1).Insert()

# Row 6 ("trialIntentTypes") already carries the label/value formatting
# (right-aligned bold label in col A, left-aligned value spanning B:E)
# that the two new rows should use, so copy its formats down into the
# freshly inserted rows 1 and 2.
$ws.Range("A6:E6").Copy()
$ws.Range("A1:E1").PasteSpecial(-4122)
$ws.Range("A2:E2").PasteSpecial(-4122)
$ws.Application.CutCopyMode = $false

# Fill in the new label/value pairs.
$ws.Range("B2").Value = "The main design for the study"
$ws.Range("B1").Value = "Study Design 1"
$ws.Range("A1").Value = "studyDesignName"
$ws.Range("A2").Value = "studyDesignDescription"

# Re-merge the value cells across B:E like the other rows in this sheet.
$ws.Range("B1:E1").Merge()
$ws.Range("B2:E2").Merge()

# This sheet becomes the active tab, with A3 selected.
$ws.Activate() | Out-Null
$ws.Range("A3").Select() | Out-Null
